$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header summary values
$ws.Range("E11").Value = 2515387
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 80

# Rows 16-64: IBETH JOHANA DE HORTA HERNANDEZ, periods 1901..2301 (monthly, ascending)
$periods1 = @("1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112","2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212","2301")

$r = 16
foreach ($p in $periods1) {
    $ws.Cells.Item($r, 3).Value = "1047365661"
    $ws.Cells.Item($r, 4).Value = "IBETH JOHANA DE HORTA HERNANDEZ"
    $ws.Cells.Item($r, 5).Value = $p
    $ws.Cells.Item($r, 6).Value = 31249
    $ws.Cells.Item($r, 7).Value = 781242
    $r = $r + 1
}

# Row 65: CELIA CRUZ MARTINEZ CORREA, period 2301
$ws.Cells.Item(65, 3).Value = "26139640"
$ws.Cells.Item(65, 4).Value = "CELIA CRUZ MARTINEZ CORREA"
$ws.Cells.Item(65, 5).Value = "2301"
$ws.Cells.Item(65, 6).Value = 15467
$ws.Cells.Item(65, 7).Value = 1423500

# Rows 66-94: IBETH JOHANA DE HORTA HERNANDEZ, periods 2302..2506 (monthly, ascending)
$periods2 = @("2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312","2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412","2501","2502","2503","2504","2505","2506")

$r = 66
foreach ($p in $periods2) {
    $ws.Cells.Item($r, 3).Value = "1047365661"
    $ws.Cells.Item($r, 4).Value = "IBETH JOHANA DE HORTA HERNANDEZ"
    $ws.Cells.Item($r, 5).Value = $p
    $ws.Cells.Item($r, 6).Value = 31249
    $ws.Cells.Item($r, 7).Value = 781242
    $r = $r + 1
}

# Row 95: IBETH JOHANA DE HORTA HERNANDEZ, period 2507
$ws.Cells.Item(95, 3).Value = "1047365661"
$ws.Cells.Item(95, 4).Value = "IBETH JOHANA DE HORTA HERNANDEZ"
$ws.Cells.Item(95, 5).Value = "2507"
$ws.Cells.Item(95, 6).Value = 31249
$ws.Cells.Item(95, 7).Value = 781242

# Row 96: IBETH JOHANA DE HORTA HERNANDEZ, period 2508
$ws.Cells.Item(96, 3).Value = "1047365661"
$ws.Cells.Item(96, 4).Value = "IBETH JOHANA DE HORTA HERNANDEZ"
$ws.Cells.Item(96, 5).Value = "2508"
$ws.Cells.Item(96, 6).Value = 31249
$ws.Cells.Item(96, 7).Value = 781242

Write-Host "Applied data updates"
